$d = $word.ActiveDocument

# 1) Remove the stale "_GoBack" bookmark left over from the previous edit
#    location (end of the MOI/aerodynamics paragraph).
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 2) Extend the "IMU ... 1 - 15kg 10 - 200 W" line with the new
#    estimate/commentary text, and drop a fresh "_GoBack" bookmark at the
#    point where the author's cursor ended up mid-sentence (right after
#    "smaller.").
$rng = $d.Content
[void]$rng.Find.Execute(" 1 – 15kg 10 – 200 W", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

$insertStart = $rng.Start
$addition = " (estimated 50W b/c smaller."
$rng.InsertAfter($addition + " gut feeling)")

$bmPos = $insertStart + $addition.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
